# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted for the Pomelo subset
# ("Vega Modelo de Temuco"). The new record belongs at row 108 (by date
# order within the sheet's existing layout), which pushes the previously
# existing rows 108-140 down to 109-141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 108; rows 108:140 shift down to 109:141,
# carrying their original values (and formatting) with them.
$ws.Rows("108:108").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A108").Value = 10
$ws.Range("B108").Value = "Vega Modelo de Temuco"
$ws.Range("C108").Value = "La Araucanía"
$ws.Range("D108").Value2 = 44463
$ws.Range("E108").Value = 9
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100102
$ws.Range("H108").Value = "Cítricos"
$ws.Range("I108").Value = 100102006
$ws.Range("J108").Value = "Pomelo"
$ws.Range("K108").Value = "Start Ruby"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 80
$ws.Range("N108").Value = 10000
$ws.Range("O108").Value = 10000
$ws.Range("P108").Value = 10000
$ws.Range("Q108").Value = "$/bandeja 15 kilos granel"
$ws.Range("R108").Value = "Región de O'Higgins"
$ws.Range("S108").Value = 667
$ws.Range("T108").Value = 15

# Keep the date column's number format consistent with the rest of
# column D (inherited from the row above on insert, set explicitly too).
$ws.Range("D108").NumberFormat = $ws.Range("D107").NumberFormat
